# Update column G ("K") values on Sheet1 to reflect corrected strikeout
# counts (regenerated from "K" instead of "Strike#").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 1
    3  = 4
    4  = 1
    5  = 7
    6  = 10
    7  = 7
    8  = 6
    9  = 8
    10 = 3
    11 = 6
    12 = 7
    13 = 4
    14 = 7
    15 = 5
    16 = 6
    17 = 6
    18 = 9
    19 = 8
    20 = 6
    21 = 10
    22 = 7
    23 = 9
    24 = 8
    25 = 8
    26 = 4
    27 = 5
    28 = 3
    29 = 5
    30 = 5
    31 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
